$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename label in A6 (shared string "DN_SRD" -> "DL_SRD")
$ws.Range("A6").Value = "DL_SRD"

# Row 3 precision refresh
$ws.Range("C3").Value = 0.01478140644530066
$ws.Range("D3").Value = 0.003458865534386073
$ws.Range("E3").Value = 0.030842967635579
$ws.Range("G3").Value = 0.01859583913606864
$ws.Range("H3").Value = 0.004585074682976598
$ws.Range("I3").Value = 0.05432804585156271

# Row 6 precision refresh
$ws.Range("B6").Value = 0.9986199039263873
$ws.Range("C6").Value = 0.03360588346129266
$ws.Range("D6").Value = 0.006771228424570977
$ws.Range("E6").Value = 0.04622462311226979
$ws.Range("F6").Value = 0.9886907418091986
$ws.Range("G6").Value = 0.1014620109562562
$ws.Range("H6").Value = 0.0226317306439741
$ws.Range("I6").Value = 0.1323231339800526
